$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new column before column I (shifts old I -> J)
$ws.Columns("I").Insert([Microsoft.Office.Interop.Excel.XlInsertShiftDirection]::xlShiftToRight)

# Header for new column (match the bold header formatting used by the
# other header cells)
$ws.Range("I1").Value = "Distribution channel code"
$ws.Range("I1").Font.Bold = $true

# New column data
$ws.Range("I2").Value = "IN"
$ws.Range("I3").Value = "GO"

# Column width for new column (COM ColumnWidth is quantized to 1/6-character
# steps by the host; 21.6666... is the input that lands closest to the
# target stored width of 22.5546875). Column J keeps the width it inherited
# automatically from the shifted-right original column I.
$ws.Columns("I").ColumnWidth = 21.666666666666668

# Update selection to match target state
$ws.Range("I3").Select()
